# Apply updated "想去人数" (want-to-go count) / ticket-status figures scraped
# at commit 456a3b4 across the four sheets of the workbook.
#
# Sheet 1 = 展览 (exhibitions), Sheet 2 = 演出 (performances),
# Sheet 3 = 本地生活 (local life), Sheet 4 = 全部类型 (all types, an aggregate
# of the first three sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 1444
$ws1.Range("F3").Value  = 1419
$ws1.Range("F5").Value  = 222
$ws1.Range("F6").Value  = 675
$ws1.Range("F7").Value  = 29
$ws1.Range("F8").Value  = 617
$ws1.Range("F9").Value  = 470
$ws1.Range("F10").Value = 75
$ws1.Range("F11").Value = 1368
$ws1.Range("F12").Value = 32083
$ws1.Range("F13").Value = 6912
$ws1.Range("F14").Value = 107
$ws1.Range("F15").Value = 347
$ws1.Range("F16").Value = 565
$ws1.Range("F17").Value = 428
$ws1.Range("F19").Value = 91
$ws1.Range("F21").Value = 434
$ws1.Range("F22").Value = 94
$ws1.Range("F23").Value = 779
$ws1.Range("F25").Value = 380
$ws1.Range("F26").Value = 427
$ws1.Range("F28").Value = 187
$ws1.Range("F30").Value = 730
$ws1.Range("F34").Value = 107
$ws1.Range("F36").Value = 779
$ws1.Range("F37").Value = 285

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F2").Value  = 1149
$ws2.Range("F3").Value  = 5
$ws2.Range("F5").Value  = 148

# Row 7 ("音波狂潮II 萤光宇宙 音游嘉年华") went back on sale: ticket count
# bumped and the lowest price cell flips from the "已售罄" (sold out) text
# back to a numeric price.
$ws2.Range("F7").Value  = 4319
$ws2.Range("G7").Value  = 480

$ws2.Range("F9").Value  = 230
$ws2.Range("F15").Value = 39
$ws2.Range("F19").Value = 4285

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F2").Value = 1430

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (aggregate of the three sheets above)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value  = 1430
$ws4.Range("F4").Value  = 1149
$ws4.Range("F5").Value  = 1444
$ws4.Range("F6").Value  = 5
$ws4.Range("F7").Value  = 1419
$ws4.Range("F8").Value  = 222
$ws4.Range("F9").Value  = 675
$ws4.Range("F10").Value = 29
$ws4.Range("F11").Value = 617
$ws4.Range("F13").Value = 75
$ws4.Range("F14").Value = 1368
$ws4.Range("F15").Value = 148

# Row 17 previously mirrored the (now stale) "昨日重现" event row; the scrape
# replaced it wholesale with the "音波狂潮II" event's current data (same
# event now reflected in Sheet2 row 7 above).
#
# B17 holds a plain-text "YYYY-MM-DD" date string (like every other cell in
# column B), not a real Excel date. Force text formatting before the write so
# COM doesn't auto-coerce the literal into a date serial number, then drop
# the format override back to Normal so no stray style survives on the cell.
$ws4.Range("B17").NumberFormat = "@"
$ws4.Range("B17").Value = "2024-07-20"
$ws4.Range("B17").Style = "Normal"
$ws4.Range("C17").Value = "广州·音波狂潮II 萤光宇宙 音游嘉年华"
$ws4.Range("D17").Value = "新港东路磨碟沙大街118号自编8栋 啤厂媒棚"
$ws4.Range("E17").Value = "2024.07.20 13:30-07.21 23:30"
$ws4.Range("F17").Value = 4319
$ws4.Range("G17").Value = 480
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=86632"
$ws4.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202405/GcK1JV3B1717123497026.jpeg"

$ws4.Range("F18").Value = 230
$ws4.Range("F19").Value = 230
$ws4.Range("F22").Value = 107
$ws4.Range("F23").Value = 347
$ws4.Range("F25").Value = 565
$ws4.Range("F26").Value = 428
$ws4.Range("F28").Value = 91
$ws4.Range("F30").Value = 39
$ws4.Range("F31").Value = 434
$ws4.Range("F32").Value = 94
$ws4.Range("F33").Value = 779
$ws4.Range("F35").Value = 380
$ws4.Range("F36").Value = 427
$ws4.Range("F38").Value = 187
$ws4.Range("F40").Value = 730
$ws4.Range("F44").Value = 107
$ws4.Range("F45").Value = 779
$ws4.Range("F46").Value = 285
